# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets.
# These mirror the same underlying events, so both sheets receive matching updates,
# just at different row offsets.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAll        = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1): cell -> new value
$exhibitionUpdates = @{
    "F2"  = 251
    "F7"  = 561
    "F8"  = 73
    "F9"  = 6818
    "F11" = 98
    "F15" = 1102
    "F16" = 16231
    "F17" = 1596
    "F18" = 41
    "F22" = 11384
    "F24" = 1029
    "F25" = 4480
    "F26" = 323
}

foreach ($addr in $exhibitionUpdates.Keys) {
    $wsExhibition.Range($addr).Value = $exhibitionUpdates[$addr]
}

# Sheet "全部类型" (sheet4): cell -> new value
$allTypesUpdates = @{
    "F2"  = 251
    "F7"  = 561
    "F9"  = 73
    "F10" = 6818
    "F12" = 98
    "F17" = 1102
    "F18" = 16231
    "F19" = 1596
    "F20" = 41
    "F26" = 11384
    "F28" = 1029
    "F29" = 4480
    "F30" = 323
}

foreach ($addr in $allTypesUpdates.Keys) {
    $wsAll.Range($addr).Value = $allTypesUpdates[$addr]
}
